$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KLBAY")

# Insert two new columns before column D (D,E) to make room for the newest
# two reporting quarters; existing D:K data shifts right to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats from the (now-shifted) original D/E columns into the
# two freshly inserted columns so the new cells render like their neighbours
# (dates in row 7/38/80, right-aligned numbers elsewhere).
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null
$ws.Range("G7:G102").Copy() | Out-Null
$ws.Range("E7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the refreshed data set (new quarter in D, prior quarter in E, and
# revised figures for the columns that shifted from D:H into F:J).
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("D8").Value = 714100
$ws.Range("E8").Value = 719600
$ws.Range("F8").Value = 573100
$ws.Range("G8").Value = 561300
$ws.Range("H8").Value = 589100
$ws.Range("I8").Value = 570300
$ws.Range("J8").Value = 508700
$ws.Range("D9").Value = 441100
$ws.Range("E9").Value = 413500
$ws.Range("F9").Value = 372100
$ws.Range("G9").Value = 399400
$ws.Range("H9").Value = 386900
$ws.Range("I9").Value = 423700
$ws.Range("J9").Value = 445600
$ws.Range("D10").Value = 273000
$ws.Range("E10").Value = 306100
$ws.Range("F10").Value = 201000
$ws.Range("G10").Value = 161900
$ws.Range("H10").Value = 202300
$ws.Range("I10").Value = 146700
$ws.Range("J10").Value = 63100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("D14").Value = -96900
$ws.Range("E14").Value = -15900
$ws.Range("F14").Value = -17700
$ws.Range("G14").Value = -30500
$ws.Range("H14").Value = -15200
$ws.Range("I14").Value = -37200
$ws.Range("J14").Value = -26100
$ws.Range("D15").Value = 2500
$ws.Range("E15").Value = 2400
$ws.Range("F15").Value = 2300
$ws.Range("G15").Value = 2300
$ws.Range("H15").Value = 2100
$ws.Range("I15").Value = 1500
$ws.Range("J15").Value = "NA"
$ws.Range("D17").Value = 433300
$ws.Range("E17").Value = 489100
$ws.Range("F17").Value = 432200
$ws.Range("G17").Value = 448500
$ws.Range("H17").Value = 453300
$ws.Range("I17").Value = 464700
$ws.Range("J17").Value = 491000
$ws.Range("D18").Value = 280800
$ws.Range("E18").Value = 230500
$ws.Range("F18").Value = 140900
$ws.Range("G18").Value = 112800
$ws.Range("H18").Value = 135800
$ws.Range("I18").Value = 105600
$ws.Range("J18").Value = 17700
$ws.Range("D20").Value = 72200
$ws.Range("E20").Value = -114200
$ws.Range("F20").Value = -441700
$ws.Range("G20").Value = 100
$ws.Range("H20").Value = -109100
$ws.Range("I20").Value = 154400
$ws.Range("J20").Value = -97800
$ws.Range("D21").Value = 459400
$ws.Range("E21").Value = 221900
$ws.Range("F21").Value = -196900
$ws.Range("G21").Value = 225900
$ws.Range("H21").Value = 123100
$ws.Range("I21").Value = 385200
$ws.Range("J21").Value = 80400
$ws.Range("D22").Value = 77800
$ws.Range("E22").Value = 79700
$ws.Range("F22").Value = 74000
$ws.Range("G22").Value = 67400
$ws.Range("H22").Value = 68500
$ws.Range("I22").Value = 69800
$ws.Range("J22").Value = 73800
$ws.Range("D23").Value = 275200
$ws.Range("E23").Value = 36500
$ws.Range("F23").Value = -374900
$ws.Range("G23").Value = 45500
$ws.Range("H23").Value = -41700
$ws.Range("I23").Value = 190200
$ws.Range("J23").Value = -153900
$ws.Range("D24").Value = 41200
$ws.Range("E24").Value = 9900
$ws.Range("F24").Value = -130100
$ws.Range("G24").Value = 13600
$ws.Range("H24").Value = -20500
$ws.Range("I24").Value = 90100
$ws.Range("J24").Value = -57100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("D26").Value = 234100
$ws.Range("E26").Value = 26600
$ws.Range("F26").Value = -244700
$ws.Range("G26").Value = 32000
$ws.Range("H26").Value = -21300
$ws.Range("I26").Value = 100100
$ws.Range("J26").Value = -96800
$ws.Range("D27").Value = 221400
$ws.Range("E27").Value = 26600
$ws.Range("F27").Value = -244700
$ws.Range("G27").Value = 32000
$ws.Range("H27").Value = -21300
$ws.Range("I27").Value = 100100
$ws.Range("J27").Value = -96800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("D32").Value = -72200
$ws.Range("E32").Value = 114200
$ws.Range("F32").Value = 441700
$ws.Range("G32").Value = -100
$ws.Range("H32").Value = 109100
$ws.Range("I32").Value = -154400
$ws.Range("J32").Value = 97800
$ws.Range("D33").Value = 221400
$ws.Range("E33").Value = 26600
$ws.Range("F33").Value = -244700
$ws.Range("G33").Value = 32000
$ws.Range("H33").Value = -21300
$ws.Range("I33").Value = 100100
$ws.Range("J33").Value = -96800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("D35").Value = 221400
$ws.Range("E35").Value = 26600
$ws.Range("F35").Value = -244700
$ws.Range("G35").Value = 32000
$ws.Range("H35").Value = -21300
$ws.Range("I35").Value = 100100
$ws.Range("J35").Value = -96800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("D41").Value = 1470000
$ws.Range("E41").Value = 1549400
$ws.Range("F41").Value = 1442300
$ws.Range("G41").Value = 1366400
$ws.Range("H41").Value = 1801900
$ws.Range("I41").Value = 1812300
$ws.Range("J41").Value = 1625400
$ws.Range("D42").Value = 336700
$ws.Range("E42").Value = 328500
$ws.Range("F42").Value = 325500
$ws.Range("G42").Value = 323300
$ws.Range("H42").Value = 318700
$ws.Range("I42").Value = 162100
$ws.Range("J42").Value = 157600
$ws.Range("D43").Value = 592400
$ws.Range("E43").Value = 576600
$ws.Range("F43").Value = 535200
$ws.Range("G43").Value = 554200
$ws.Range("H43").Value = 595100
$ws.Range("I43").Value = 613600
$ws.Range("J43").Value = 548500
$ws.Range("D44").Value = 309300
$ws.Range("E44").Value = 289300
$ws.Range("F44").Value = 266900
$ws.Range("G44").Value = 252300
$ws.Range("H44").Value = 239200
$ws.Range("I44").Value = 232400
$ws.Range("J44").Value = 239100
$ws.Range("D45").Value = 76300
$ws.Range("E45").Value = 64800
$ws.Range("F45").Value = 68400
$ws.Range("G45").Value = 63600
$ws.Range("H45").Value = 71200
$ws.Range("I45").Value = 70900
$ws.Range("J45").Value = 59200
$ws.Range("D46").Value = 2784800
$ws.Range("E46").Value = 2808600
$ws.Range("F46").Value = 2638400
$ws.Range("G46").Value = 2559700
$ws.Range("H46").Value = 3026200
$ws.Range("I46").Value = 2891300
$ws.Range("J46").Value = 2629700
$ws.Range("D47").Value = 370800
$ws.Range("E47").Value = 389400
$ws.Range("F47").Value = 406100
$ws.Range("G47").Value = 363900
$ws.Range("H47").Value = 374100
$ws.Range("I47").Value = 413400
$ws.Range("J47").Value = 426400
$ws.Range("D48").Value = 4320700
$ws.Range("E48").Value = 4200900
$ws.Range("F48").Value = 4218000
$ws.Range("G48").Value = 4250300
$ws.Range("H48").Value = 4299200
$ws.Range("I48").Value = 4294900
$ws.Range("J48").Value = 4340000
$ws.Range("D49").Value = 21800
$ws.Range("E49").Value = 22400
$ws.Range("F49").Value = 22800
$ws.Range("G49").Value = 22800
$ws.Range("H49").Value = 23100
$ws.Range("I49").Value = 23500
$ws.Range("J49").Value = 23900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("D52").Value = 99300
$ws.Range("E52").Value = 106800
$ws.Range("F52").Value = 107900
$ws.Range("G52").Value = 113800
$ws.Range("H52").Value = 109600
$ws.Range("I52").Value = 109700
$ws.Range("J52").Value = 112100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("D54").Value = 7597500
$ws.Range("E54").Value = 7528100
$ws.Range("F54").Value = 7393200
$ws.Range("G54").Value = 7310400
$ws.Range("H54").Value = 7832300
$ws.Range("I54").Value = 7732900
$ws.Range("J54").Value = 7532100
$ws.Range("D57").Value = 231700
$ws.Range("E57").Value = 212500
$ws.Range("F57").Value = 206900
$ws.Range("G57").Value = 158300
$ws.Range("H57").Value = 183000
$ws.Range("I57").Value = 148400
$ws.Range("J57").Value = 157800
$ws.Range("D58").Value = 506500
$ws.Range("E58").Value = 582900
$ws.Range("F58").Value = 563600
$ws.Range("G58").Value = 457800
$ws.Range("H58").Value = 633200
$ws.Range("I58").Value = 600000
$ws.Range("J58").Value = 599800
$ws.Range("D59").Value = 212600
$ws.Range("E59").Value = 159000
$ws.Range("F59").Value = 139800
$ws.Range("G59").Value = 131500
$ws.Range("H59").Value = 144600
$ws.Range("I59").Value = 146100
$ws.Range("J59").Value = 154400
$ws.Range("D60").Value = 950700
$ws.Range("E60").Value = 954400
$ws.Range("F60").Value = 910400
$ws.Range("G60").Value = 747600
$ws.Range("H60").Value = 960700
$ws.Range("I60").Value = 894500
$ws.Range("J60").Value = 912000
$ws.Range("D61").Value = 4479000
$ws.Range("E61").Value = 4580900
$ws.Range("F61").Value = 4433600
$ws.Range("G61").Value = 4079600
$ws.Range("H61").Value = 4378800
$ws.Range("I61").Value = 4232200
$ws.Range("J61").Value = 4195100
$ws.Range("D62").Value = 492800
$ws.Range("E62").Value = 449400
$ws.Range("F62").Value = 487200
$ws.Range("G62").Value = 637000
$ws.Range("H62").Value = 638000
$ws.Range("I62").Value = 684000
$ws.Range("J62").Value = 597200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("D66").Value = 5984300
$ws.Range("E66").Value = 5984700
$ws.Range("F66").Value = 5831200
$ws.Range("G66").Value = 5464200
$ws.Range("H66").Value = 5977600
$ws.Range("I66").Value = 5810700
$ws.Range("J66").Value = 5704300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("D72").Value = 355600
$ws.Range("E72").Value = 285100
$ws.Range("F72").Value = 303900
$ws.Range("G72").Value = 587600
$ws.Range("H72").Value = 996500
$ws.Range("I72").Value = 1089100
$ws.Range("J72").Value = 994200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("D76").Value = 1613200
$ws.Range("E76").Value = 1543400
$ws.Range("F76").Value = 1562000
$ws.Range("G76").Value = 1846200
$ws.Range("H76").Value = 1854700
$ws.Range("I76").Value = 1922200
$ws.Range("J76").Value = 1827700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("D81").Value = 221400
$ws.Range("E81").Value = 26600
$ws.Range("F81").Value = -244700
$ws.Range("G81").Value = 32000
$ws.Range("H81").Value = -21300
$ws.Range("I81").Value = 100100
$ws.Range("J81").Value = -96800
$ws.Range("D83").Value = 106400
$ws.Range("E83").Value = 105700
$ws.Range("F83").Value = 104000
$ws.Range("G83").Value = 113000
$ws.Range("H83").Value = 96400
$ws.Range("I83").Value = 125200
$ws.Range("J83").Value = 160400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("D89").Value = 169800
$ws.Range("E89").Value = 174200
$ws.Range("F89").Value = 239600
$ws.Range("G89").Value = 130900
$ws.Range("H89").Value = 86200
$ws.Range("I89").Value = 82700
$ws.Range("J89").Value = 135400
$ws.Range("D91").Value = -58300
$ws.Range("E91").Value = -23300
$ws.Range("F91").Value = -21100
$ws.Range("G91").Value = -15700
$ws.Range("H91").Value = -18600
$ws.Range("I91").Value = -14400
$ws.Range("J91").Value = -16600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("D94").Value = -105200
$ws.Range("E94").Value = -62700
$ws.Range("F94").Value = -45200
$ws.Range("G94").Value = -55900
$ws.Range("H94").Value = -60800
$ws.Range("I94").Value = -54100
$ws.Range("J94").Value = -36700
$ws.Range("D96").Value = -400
$ws.Range("E96").Value = 80800
$ws.Range("F96").Value = -39000
$ws.Range("G96").Value = -43800
$ws.Range("H96").Value = 81000
$ws.Range("I96").Value = -45400
$ws.Range("J96").Value = -12800
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("D100").Value = -144000
$ws.Range("E100").Value = -4400
$ws.Range("F100").Value = -118400
$ws.Range("G100").Value = -510600
$ws.Range("H100").Value = -35800
$ws.Range("I100").Value = 158300
$ws.Range("J100").Value = -178900
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("D102").Value = -79400
$ws.Range("E102").Value = 107100
$ws.Range("F102").Value = 75900
$ws.Range("G102").Value = -435600
$ws.Range("H102").Value = -10300
$ws.Range("I102").Value = 186900
$ws.Range("J102").Value = -80200
